$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 8340
    $ws.Range("F3").Value = 7772
    $ws.Range("F10").Value = 163
    $ws.Range("F12").Value = 706
    $ws.Range("F14").Value = 1327
    $ws.Range("F15").Value = 61
    $ws.Range("F17").Value = 11
    $ws.Range("F19").Value = 121
}

$wb.Save()
